# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Fri Sep 20 13:56:20 UTC 2024 with GitHub Actions"
#
# Column D holds price strings and column E holds padded "  +x.xx%  " strings;
# both are stored as plain text (t="inlineStr") in the source workbook, not
# numbers. Assigning a numeric-looking string straight to .Value lets the COM
# layer coerce it into a real number (e.g. "1.00" -> 1), which would corrupt the
# cell type. Forcing NumberFormat "@" (Text) before the write, then clearing the
# format again right after, keeps the write as literal text while leaving the
# cell style back at its original (default/unstyled) state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

Set-TextCell "D2" "62.879.52"
$ws.Range("E2").Value = "  +0.11%  "
Set-TextCell "D3" "2.528.84"
$ws.Range("E3").Value = "  +4.29%  "
$ws.Range("E4").Value = "  -0.09%  "
Set-TextCell "D5" "567.42"
$ws.Range("E5").Value = "  +0.84%  "
Set-TextCell "D6" "147.54"
$ws.Range("E6").Value = "  +5.18%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  -0.38%  "
Set-TextCell "D9" "2.527.74"
$ws.Range("E9").Value = "  +4.31%  "
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  -1.52%  "
$ws.Range("E12").Value = "  +1.14%  "
$ws.Range("E13").Value = "  +1.51%  "
Set-TextCell "D14" "27.69"
$ws.Range("E14").Value = "  +5.49%  "
Set-TextCell "D15" "2.978.89"
$ws.Range("E15").Value = "  +4.15%  "
Set-TextCell "D16" "62.836.51"
$ws.Range("E16").Value = "  +0.23%  "
$ws.Range("E17").Value = "  -0.02%  "
Set-TextCell "D18" "2.537.90"
$ws.Range("E18").Value = "  +4.44%  "
$ws.Range("E19").Value = "  +3.32%  "
Set-TextCell "D20" "335.40"
$ws.Range("E20").Value = "  -1.12%  "
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("E22").Value = "  -0.33%  "
$ws.Range("E23").Value = "  -0.02%  "
Set-TextCell "D24" "65.57"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("B26").Value = "SuiNetwork"
$ws.Range("C26").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D26" "1.52"
$ws.Range("E26").Value = "  +13.29%  "
$ws.Range("B27").Value = "Fetch.AI"
$ws.Range("C27").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextCell "D27" "1.57"
$ws.Range("E27").Value = "  +3.12%  "
Set-TextCell "D28" "0.999"
$ws.Range("E28").Value = "  -0.15%  "
Set-TextCell "D29" "8.34"
$ws.Range("E29").Value = "  +2.30%  "
Set-TextCell "D30" "7.24"
$ws.Range("E30").Value = "  +11.36%  "
Set-TextCell "D31" "0.0₃0811"
$ws.Range("E31").Value = "  +2.35%  "
Set-TextCell "D32" "1.84"
$ws.Range("E32").Value = "  +1.13%  "
Set-TextCell "D33" "177.75"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("E34").Value = "  +7.70%  "
Set-TextCell "D35" "411.37"
$ws.Range("E35").Value = "  +11.25%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("E40").Value = "  +2.11%  "
Set-TextCell "D41" "1.00"
$ws.Range("E41").Value = "  +0.14%  "
Set-TextCell "D42" "39.23"
$ws.Range("E42").Value = "  -1.52%  "
Set-TextCell "D43" "151.76"
$ws.Range("E43").Value = "  +4.11%  "
Set-TextCell "D44" "3.75"
$ws.Range("E44").Value = "  +1.95%  "
Set-TextCell "D45" "20.66"
Set-TextCell "D46" "0.603"
$ws.Range("E46").Value = "  +2.40%  "
$ws.Range("E47").Value = "  +1.55%  "
$ws.Range("E48").Value = "  +0.38%  "
Set-TextCell "D49" "0.0237"
$ws.Range("E49").Value = "  +6.32%  "
Set-TextCell "D50" "18.30"
$ws.Range("E50").Value = "  +2.80%  "
Set-TextCell "D51" "1.78"
$ws.Range("E51").Value = "  +3.39%  "
